# Update stats for 2026-01 (row 26) in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B26").Value = 6488
$ws.Range("C26").Value = 1012
$ws.Range("D26").Value = 6042745
$ws.Range("E26").Value = 931.3725339087546
$ws.Range("F26").Value = 9.650160554334963
$ws.Range("G26").Value = 7.430997876857748
$ws.Range("H26").Value = 25.84036820325566
